$d = $word.ActiveDocument

# The target is the last paragraph in the document - an empty bullet
# list item ("w14:paraId=2B38A461") that gets five new runs (with
# spell-check proofErr markers around the English function call)
# describing use of get_the_date('j F Y').
$p = $d.Paragraphs.Last

$r = $p.Range

# Pull the canonical OOXML for this range so we can splice in the new
# runs with byte-exact formatting (rFonts/sz/szCs/rtl/lang) instead of
# relying on Range.InsertAfter, which only ever records "diffs" against
# inherited paragraph-mark formatting and would silently drop the
# explicit w:cs="Shabnam" etc. that the target XML requires.
$xml = $r.WordOpenXML

$oldPara = '<w:p w14:paraId="2B38A461" w14:textId="77777777" w:rsidR="00B96B08" w:rsidRPr="00BA6B77" w:rsidRDefault="00B96B08" w:rsidP="00B96B08"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:bidi/><w:rPr><w:rFonts w:ascii="Shabnam" w:hAnsi="Shabnam" w:cs="Shabnam"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr></w:p>'

$newRuns = '<w:r><w:rPr><w:rFonts w:ascii="Shabnam" w:hAnsi="Shabnam" w:cs="Shabnam" w:hint="cs"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">از </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Shabnam" w:hAnsi="Shabnam" w:cs="Shabnam"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>get_the_date</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Shabnam" w:hAnsi="Shabnam" w:cs="Shabnam"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>(''j F Y'')</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Shabnam" w:hAnsi="Shabnam" w:cs="Shabnam"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Shabnam" w:hAnsi="Shabnam" w:cs="Shabnam" w:hint="cs"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> برای درست کردن ترتیب تاریخ ها توی بلاگ استفاده کردم</w:t></w:r>'

$newPara = '<w:p w14:paraId="2B38A461" w14:textId="77777777" w:rsidR="00B96B08" w:rsidRPr="00BA6B77" w:rsidRDefault="00B96B08" w:rsidP="00B96B08"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:bidi/><w:rPr><w:rFonts w:ascii="Shabnam" w:hAnsi="Shabnam" w:cs="Shabnam"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr>' + $newRuns + '</w:p>'

if (-not $xml.Contains($oldPara)) {
    throw "Could not locate target empty paragraph in WordOpenXML"
}

$newXml = $xml.Replace($oldPara, $newPara)

$r.InsertXML($newXml) | Out-Null

Write-Output "Inserted get_the_date() note into final paragraph."
